$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the header row's formatting (bold font, border, alignment) by
# copying it to a scratch area before the sheet contents are cleared.
$ws.Range("A1:E1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Clear out the old table (rows 1-9) so the shared-string pool is rebuilt
# cleanly and only contains strings actually used by the new table.
$ws.Range("A1:E9").Clear()

# Row 1 - header (unchanged content)
$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "Site"
$ws.Range("C1").Value = "Exposure period"
$ws.Range("D1").Value = "Notes"
$ws.Range("E1").Value = "Exist"

# Re-apply the header formatting that was saved off above, then remove the
# scratch copy.
$ws.Range("G1:K1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1:K1").Clear()
$excel.CutCopyMode = 0

# Row 2
$ws.Range("A2").Value = "30/12/20 10:45am- 12:15pm"
$ws.Range("E2").Value = "old"

# Row 3 - Cheltenham, Angus and Cootes Jeweller (old address wording)
$ws.Range("A3").Value = "Cheltenham"
$ws.Range("B3").Value = "Angus and Cootes Jeweller  Southland Shopping Centre, 2096/1239 Nepean Hwy, Cheltenham VIC 3192"
$ws.Range("C3").Value = "28/12/2020 2:30pm-2:50pm"
$ws.Range("D3").Value = "Case shopped in store"
$ws.Range("E3").Value = "old"

# Row 4 - Cheltenham, Angus and Cootes Jeweller (new address wording)
$ws.Range("A4").Value = "Cheltenham"
$ws.Range("B4").Value = "Angus and Cootes Jeweller  Southland Shopping Centre, Shop 2096/1239, Nepean Hwy, Cheltenham VIC 3192"
$ws.Range("C4").Value = "28/12/2020 2:30pm-2:50pm"
$ws.Range("D4").Value = "Case shopped in store"
$ws.Range("E4").Value = "new"

# Row 5 - Moorabbin, COSTCO (old exposure period with typo)
$ws.Range("A5").Value = "Moorabbin"
$ws.Range("B5").Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Range("C5").Value = "30/12/20 10:45am- 12:15pm"
$ws.Range("D5").Value = "Case shopped in store"
$ws.Range("E5").Value = "old"

# Row 6 - Moorabbin, COSTCO (new exposure period, typo fixed)
$ws.Range("A6").Value = "Moorabbin"
$ws.Range("B6").Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Range("C6").Value = "30/12/20 10:45am-12:15pm"
$ws.Range("D6").Value = "Case shopped in store"
$ws.Range("E6").Value = "new"

# The table now only spans rows 1-6; the longer replacement text in columns
# B and C needs wider "best fit" columns to display without truncation.
$ws.Columns("B").ColumnWidth = 87.9296875
$ws.Columns("C").ColumnWidth = 24
